$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting rows 58:155 down to 59:156.
$ws.Rows("58:58").Insert()

# Populate the newly inserted row 58 with the new data point.
$ws.Range("A58").Value2 = 10
$ws.Range("B58").Value2 = "Vega Modelo de Temuco"
$ws.Range("C58").Value2 = "La Araucanía"
$ws.Range("D58").Value2 = 44519
$ws.Range("E58").Value2 = 9
$ws.Range("F58").Value2 = 100112013
$ws.Range("G58").Value2 = "Alcachofa"
$ws.Range("H58").Value2 = "Madrigal"
$ws.Range("I58").Value2 = "Primera"
$ws.Range("J58").Value2 = 65
$ws.Range("K58").Value2 = 12000
$ws.Range("L58").Value2 = 12000
$ws.Range("M58").Value2 = 12000
$ws.Range("N58").Value2 = "$/caja 40 unidades"
$ws.Range("O58").Value2 = "Región del Maule"
$ws.Range("P58").Value2 = 300
$ws.Range("Q58").Value2 = 40
$ws.Range("R58").Value2 = "Hortaliza"
